# Updated cryptos list with GitHub Actions - apply latest price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force the cell to keep a purely-numeric-looking string as text,
    # matching the original inline-string storage instead of letting
    # Excel auto-convert it to a floating point number.
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.143.47"
$ws.Range("E2").Value = "  -2.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.133.28"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "593.39"
$ws.Range("E5").Value = "  -2.58%  "

# Row 6 - Solana
Set-TextValue "D6" "136.50"
$ws.Range("E6").Value = "  -5.13%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.122.48"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -2.43%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.32%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.22"
$ws.Range("E11").Value = "  -3.19%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.458"
$ws.Range("E12").Value = "  -3.15%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.13%  "

# Row 14 - Avalanche
Set-TextValue "D14" "34.25"
$ws.Range("E14").Value = "  -3.42%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.638.01"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +2.47%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.122.51"
$ws.Range("E17").Value = "  -2.04%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.126.38"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.70"
$ws.Range("E19").Value = "  -2.60%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "474.74"
$ws.Range("E20").Value = "  -0.45%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.26"
$ws.Range("E21").Value = "  -3.86%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.700"
$ws.Range("E22").Value = "  -2.79%  "

# Row 23 - Uniswap
Set-TextValue "D23" "7.74"
$ws.Range("E23").Value = "  -0.84%  "

# Row 24 - Litecoin
Set-TextValue "D24" "86.87"

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "13.04"
$ws.Range("E25").Value = "  -4.33%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -1.87%  "

# Row 28 - NEARProtocol
Set-TextValue "D28" "7.18"
$ws.Range("E28").Value = "  -2.91%  "

# Row 29 - RenderToken
Set-TextValue "D29" "7.94"
$ws.Range("E29").Value = "  -6.25%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -0.46%  "

# Row 31 - FirstDigitalUSD
$ws.Range("E31").Value = "  +0.02%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "26.74"
$ws.Range("E32").Value = "  +0.22%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -7.06%  "

# Row 34 - Stacks
$ws.Range("E34").Value = "  -4.11%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -3.26%  "

# Row 36 - Filecoin
Set-TextValue "D36" "5.83"
$ws.Range("E36").Value = "  -2.59%  "

# Row 37 - OKB
Set-TextValue "D37" "52.05"
$ws.Range("E37").Value = "  -1.22%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -5.09%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -1.84%  "

# Row 40 - Bittensor
Set-TextValue "D40" "421.89"
$ws.Range("E40").Value = "  -6.40%  "

# Row 41 - Cosmos
$ws.Range("E41").Value = "  -1.03%  "

# Row 42 - dogwifhat
$ws.Range("E42").Value = "  -9.78%  "

# Row 43 - was Maker, now Kaspa (rows swapped with old row 44)
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D43" "0.114"
$ws.Range("E43").Value = "  -3.80%  "

# Row 44 - was Kaspa, now Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.887.45"
$ws.Range("E44").Value = "  +0.36%  "

# Row 45 - TheGraph
Set-TextValue "D45" "0.263"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46 - was Fetch.AI, now USDe (rows swapped with old row 47)
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47 - was USDe, now Fetch.AI
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D47" "2.13"
$ws.Range("E47").Value = "  -4.42%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "25.81"
$ws.Range("E48").Value = "  -2.46%  "

# Row 49 - was Stellar, now ThetaToken (rows swapped with old row 50)
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D49" "2.29"
$ws.Range("E49").Value = "  -5.74%  "

# Row 50 - was ThetaToken, now Stellar
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.113"
$ws.Range("E50").Value = "  -0.95%  "

# Row 51 - Monero
Set-TextValue "D51" "120.13"
$ws.Range("E51").Value = "  -0.85%  "
